$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The QA tenant credentials (Tenant/Username/Password) in row 2 are no
# longer relevant -- clear them, leaving just the Testing URL cell.
$ws.Range("B2:D2").ClearContents()

# Point the Testing URL cell at the new environment and make it a real
# hyperlink (re-inserting it applies Excel's built-in "Hyperlink" cell
# style: underlined, theme-colored font).
$ws.Range("A2").Value = "https://replace.rmx.rentmanager.qa/"
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://replace.rmx.rentmanager.qa/")

$ws.Range("A2").Select() | Out-Null
